$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 additions (X3, Y3) ---
$ws.Range("X3").Value = -0.34999899999999684
$ws.Range("Y3").Value = "Down"

# --- Row 4 (new) ---
$ws.Range("A4").Value = 42649.612187500003
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = "Buy"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = "Random"
$ws.Range("Q4").Value = 35.483823948801813
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0.0965
$ws.Range("T4").Value = 0.0269
$ws.Range("U4").Value = 4.82
$ws.Range("V4").Value = 2.2799999999999998
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = -0.34999899999999684
$ws.Range("Y4").Value = "Down"

# --- Row 5 (new) ---
$ws.Range("A5").Value = 42649.635567129626
$ws.Range("B5").Value = 8
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = "Random"
$ws.Range("Q5").Value = 35.483823948801813
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0.0965
$ws.Range("T5").Value = 0.0269
$ws.Range("U5").Value = 4.82
$ws.Range("V5").Value = 2.2799999999999998
$ws.Range("W5").Value = 0

# --- Apply number formats matching rows above (date format for A, percent for S/T), preserving shared style usage ---
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)

$ws.Range("S3:T3").Copy()
$ws.Range("S4:T5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

Write-Host "done"
